$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Update the TestCases value in B2 from 42 to 40 (keep it as text, matching shared string type)
$ws.Range("B2").Value = "40"

# Update the selected cell on the sheet to C2 (was E2)
$ws.Range("C2").Select()
